$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 85; this pushes the existing rows 85..155
# down to 86..156 (matching the diff: every D/J/K/L/M/P value below row 85
# is inherited from the row that used to be one above it, and the former
# row 155 becomes the new row 156).
$ws.Rows("85:85").Insert()

# Populate the newly inserted row 85 with the new weekly price record.
$ws.Range("A85").Value = 4
$ws.Range("B85").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C85").Value = "Los Lagos"
$ws.Range("D85").Value = 44447
$ws.Range("E85").Value = 10
$ws.Range("F85").Value = 100112037
$ws.Range("G85").Value = "Cebollín"
$ws.Range("H85").Value = "Sin especificar"
$ws.Range("I85").Value = "Primera"
$ws.Range("J85").Value = 36
$ws.Range("K85").Value = 6000
$ws.Range("L85").Value = 6000
$ws.Range("M85").Value = 6000
$ws.Range("N85").Value = "$/paquete 36 unidades"
$ws.Range("O85").Value = "Región Metropolitana"
$ws.Range("P85").Value = 167
$ws.Range("Q85").Value = 36
$ws.Range("R85").Value = "Hortaliza"

# Keep the date column's number format consistent with the rest of column D.
$ws.Range("D85").NumberFormat = $ws.Range("D86").NumberFormat
